# Add most recent YPE data from Google Drive
#
# This applies a batch of species/association-code relabels coming from an
# updated Google Sheets export of the YPE_74 understory plot data:
#
#   forbe_unknown -> forb_YPE74      (assoc2 column; also subsumes the
#                                      separate "forb_unknown" label used
#                                      in later rows - both collapse to the
#                                      same new label)
#   agrostis      -> Agrostis_sp     (assoc3 column)
#   GATR          -> Galium_1        (assoc4 column)
#   LULE          -> Lupinus_1       (assoc8 column)
#   senecio       -> Asteraceae_1    (assoc9 column)
#   IRTE          -> Iris_1          (assoc11 column)
#
# Everything else in the sheet (plot metadata, dates, species1-10 codes
# such as rock/litter/wood/dirt/ARVI/QUKE/PIPO, remaining assoc columns)
# is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$used = $ws.UsedRange

$xlWhole = 1

# "forb_unknown" and "forbe_unknown" are distinct labels in the original
# data that both become "forb_YPE74"; replacing both (in either order)
# merges them into the single new label.
$used.Replace("forbe_unknown", "forb_YPE74", $xlWhole)
$used.Replace("forb_unknown", "forb_YPE74", $xlWhole)

$used.Replace("agrostis", "Agrostis_sp", $xlWhole)
$used.Replace("GATR", "Galium_1", $xlWhole)
$used.Replace("LULE", "Lupinus_1", $xlWhole)
$used.Replace("senecio", "Asteraceae_1", $xlWhole)
$used.Replace("IRTE", "Iris_1", $xlWhole)
